# zeroPadding no id, caso FlowMap possua scenarioPrefix
#
# Rewrites the ID column (A) of the exported-scenarios sheet from plain
# sequential numbers ("1", "2", ... "192") to a zero-padded id carrying the
# FlowMap's scenario prefix ("02.001", "02.002", ... "02.192").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$prefix = "02"

# Find the extent of the used range (header row + one row per scenario).
$lastRow = $ws.UsedRange.Rows.Count

# Mark the ID column as Text up-front (one Range.NumberFormat call -> one
# shared style for every data row) so the zero-padded strings we assign
# below aren't re-interpreted by Excel as numbers (which would silently
# drop the leading zero, e.g. "02.001" -> "2.001").
$idRange = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item($lastRow, 1))
$idRange.NumberFormat = "@"

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $oldId = [int]$cell.Value2
    $newId = "{0}.{1:D3}" -f $prefix, $oldId
    $cell.Value = $newId
}
